$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-25, replacing the previous
# Strike# derived values with the regenerated K values.
$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 2
    24 = 0
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
